# Update "想去人数" (interest count) values in the F column across sheets,
# reflecting newly scraped counts from the site.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 1750
$ws1.Range("F8").Value  = 324
$ws1.Range("F9").Value  = 297
$ws1.Range("F10").Value = 1680
$ws1.Range("F12").Value = 1389
$ws1.Range("F15").Value = 655
$ws1.Range("F16").Value = 12625
$ws1.Range("F17").Value = 12653
$ws1.Range("F23").Value = 496
$ws1.Range("F24").Value = 1977

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 126
$ws2.Range("F9").Value = 49

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value  = 1750
$ws4.Range("F12").Value = 324
$ws4.Range("F14").Value = 297
$ws4.Range("F15").Value = 1680
$ws4.Range("F17").Value = 1389
$ws4.Range("F21").Value = 655
$ws4.Range("F22").Value = 12625
$ws4.Range("F23").Value = 12653
$ws4.Range("F29").Value = 496
$ws4.Range("F32").Value = 1977
$ws4.Range("F34").Value = 126
$ws4.Range("F36").Value = 49
